$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "41.525.65"; DForceText = $false; E = "  +1.27%  " },
    @{ Row = 3; D = "2.477.41"; DForceText = $false; E = "  +0.88%  " },
    @{ Row = 4; E = "  -0.07%  " },
    @{ Row = 5; D = "313.25"; DForceText = $true; E = "  +0.80%  " },
    @{ Row = 6; D = "93.31"; DForceText = $true; E = "  +0.41%  " },
    @{ Row = 7; E = "  -1.21%  " },
    @{ Row = 8; E = "  -0.05%  " },
    @{ Row = 9; E = "  +2.43%  " },
    @{ Row = 10; D = "32.71"; DForceText = $true; E = "  -1.04%  " },
    @{ Row = 11; E = "  +1.46%  " },
    @{ Row = 12; E = "  +2.91%  " },
    @{ Row = 13; D = "2.858.92"; DForceText = $false; E = "  +0.94%  " },
    @{ Row = 14; D = "6.82"; DForceText = $true; E = "  -1.55%  " },
    @{ Row = 15; D = "16.02"; DForceText = $true; E = "  +8.82%  " },
    @{ Row = 16; D = "2.498.45"; DForceText = $false; E = "  +1.55%  " },
    @{ Row = 17; E = "  -1.71%  " },
    @{ Row = 18; D = "41.503.00"; DForceText = $false; E = "  +1.27%  " },
    @{ Row = 19; D = "6.40"; DForceText = $true; E = "  +2.44%  " },
    @{ Row = 20; E = "  +2.78%  " },
    @{ Row = 21; D = "71.58"; DForceText = $true; E = "  +5.54%  " },
    @{ Row = 22; D = "11.31"; DForceText = $true; E = "  +2.69%  " },
    @{ Row = 23; D = "236.40"; DForceText = $true; E = "  +0.97%  " },
    @{ Row = 24; D = "2.70"; DForceText = $true; E = "  -1.06%  " },
    @{ Row = 25; E = "  -0.42%  " },
    @{ Row = 26; E = "  +0.38%  " },
    @{ Row = 27; D = "24.81"; DForceText = $true; E = "  +4.67%  " },
    @{ Row = 28; E = "  +0.22%  " },
    @{ Row = 29; E = "  +1.32%  " },
    @{ Row = 30; D = "35.76"; DForceText = $true; E = "  +0.60%  " },
    @{ Row = 31; D = "158.28"; DForceText = $true; E = "  +4.81%  " },
    @{ Row = 32; D = "5.48"; DForceText = $true; E = "  +0.53%  " },
    @{ Row = 33; D = "2.57"; DForceText = $true; E = "  +1.78%  " },
    @{ Row = 34; E = "  +2.75%  " },
    @{ Row = 35; D = "17.50"; DForceText = $true; E = "  +5.48%  " },
    @{ Row = 36; E = "  -7.92%  " },
    @{ Row = 37; D = "2.91"; DForceText = $true; E = "  -1.56%  " },
    @{ Row = 38; E = "  +3.75%  " },
    @{ Row = 39; E = "  -0.82%  " },
    @{ Row = 40; E = "  +0.35%  " },
    @{ Row = 41; D = "4.12"; DForceText = $true; E = "  -0.53%  " },
    @{ Row = 42; E = "  -0.21%  " },
    @{ Row = 43; D = "19.53"; DForceText = $true; E = "  -1.49%  " },
    @{ Row = 44; D = "1.971.08"; DForceText = $false; E = "  +0.34%  " },
    @{ Row = 45; E = "  +0.73%  " },
    @{ Row = 46; E = "  -1.14%  " },
    @{ Row = 47; D = "9.10"; DForceText = $true; E = "  +6.39%  " },
    @{ Row = 48; D = "2.719.29"; DForceText = $false; E = "  +1.16%  " },
    @{ Row = 49; D = "98.13"; DForceText = $true; E = "  +1.95%  " },
    @{ Row = 50; D = "68.03"; DForceText = $true; E = "  -1.71%  " },
    @{ Row = 51; D = "72.35"; DForceText = $true; E = "  -1.80%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($r, 4)
        if ($u.DForceText) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
